# Updated model class and types
# Adds a new "SPHARM" model-type row (3D / framework / SPHARM / N/A /
# cell membrane, nuclear membrane) to the "v2.8.0" sheet, right below the
# existing data, and moves the active selection to follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v2.8.0")
$ws.Activate()

# New row of data (row 15) appended under the existing table contents.
$ws.Range("A15").Value = "3D"
$ws.Range("B15").Value = "framework"
$ws.Range("C15").Value = "SPHARM"
$ws.Range("D15").Value = "N/A"
$ws.Range("E15").Value = "cell membrane, nuclear membrane"

# Move the visible selection, matching the post-edit cursor position.
[void]$ws.Range("E16").Select()
